$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update selection (active cell) on the sheet view
$ws.Range("F12").Select()

# Update TSR formulas in column F (rows 4-6)
$ws.Range("F4").Formula = "= (5 * SUM(1, -E4) + 4 * SUM(1, -D4)) / 9"
$ws.Range("F5:F6").Formula = "= (5 * SUM(1, -E5) + 4 * SUM(1, -D5)) / 9"
